$d = $word.ActiveDocument

# 1. Remove the _GoBack bookmark currently sitting on the
#    "How are Anomalies and Suspicious behaviour different?" heading.
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# 2. Locate the "Problem Statement" heading paragraph and the empty
#    separator paragraph that immediately follows it.
$target = $null
$nextP = $null
$found = $false
foreach ($p in $d.Paragraphs) {
    if ($found) {
        $nextP = $p
        $found = $false
    }
    if ($p.Range.Text -like "Problem Statement*") {
        $target = $p
        $found = $true
    }
}

# 3. Insert a new empty paragraph + "Related Work" Heading 2 paragraph
#    (carrying the _GoBack bookmark) right after that separator, while
#    preserving a blank separator before the following "Organisation"
#    heading.
$rng = $nextP.Range
$rng.Collapse(1)

$xml = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>
<w:p/><w:p><w:pPr><w:pStyle w:val="Heading2"/></w:pPr><w:r><w:t>Related</w:t></w:r><w:r><w:t xml:space="preserve"> Work</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p><w:p/><w:p/>
</w:body></w:document>
</pkg:xmlData></pkg:part></pkg:package>
'@
$rng.InsertXML($xml)
